# Move the 4 "New" listings into "Previously added" (appended at the bottom)
# and replace "New" with 3 freshly scraped listings.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# Helper: force a value to be written as literal TEXT (never auto-coerced to a
# number/date by Excel) without creating any new cell style. We do this by
# writing a `="<text>"` formula and then immediately flattening it back down
# to a plain value with a self Copy / PasteSpecial(values).
function Set-TextValue {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$lastRow1 = $ws1.UsedRange.Rows.Count
$moveCount = $ws2.UsedRange.Rows.Count - 1

# --- 1. Copy the existing "New" rows (2..end) down to the bottom of
#        "Previously added", preserving values, hyperlinks and formatting ---
for ($i = 0; $i -lt $moveCount; $i++) {
    $srcRow = 2 + $i
    $dstRow = $lastRow1 + 1 + $i

    # Column A always holds the literal link URL, so it doubles as the
    # hyperlink's target address (reading back `.Address` on a hyperlink
    # that was loaded from the file - rather than added this session -
    # isn't reliable here).
    $srcLinkAddress = [string]$ws2.Cells.Item($srcRow, 1).Value2

    # Clone formatting from the last existing data row so the new rows get
    # the exact same cell styles (no new style entries are introduced).
    $ws1.Range("A" + $lastRow1 + ":F" + $lastRow1).Copy() | Out-Null
    $ws1.Range("A" + $dstRow + ":F" + $dstRow).PasteSpecial(-4122) | Out-Null

    for ($c = 1; $c -le 5; $c++) {
        $val = [string]$ws2.Cells.Item($srcRow, $c).Value2
        if ($val -ne "") {
            Set-TextValue $ws1.Cells.Item($dstRow, $c) $val
        }
    }
    $ws1.Cells.Item($dstRow, 6).Value = $ws2.Cells.Item($srcRow, 6).Value2

    if ($srcLinkAddress) {
        $ws1.Hyperlinks.Add($ws1.Cells.Item($dstRow, 1), $srcLinkAddress) | Out-Null
    }
}

# --- 2. Clear out the "New" sheet (data rows + their hyperlinks) ---
$ws2.Hyperlinks.Delete() | Out-Null
$ws2.Rows("2:" + $ws2.UsedRange.Rows.Count).Delete() | Out-Null

# --- 3. Populate "New" with the freshly scraped listings ---
$newListings = @(
    @("https://www.ss.com/msg/lv/real-estate/wood/bauska-and-reg/iecavas-nov/lclfm.html", "7 000 €", "Bauska un raj.", "4 ha.", "40460060490", 46071.69236111111),
    @("https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/kepovas-pag/mxnne.html", "65 000 €", "Krāslava un raj.", "9.50 ha.", "60800040007", 46072.018055555556),
    @("https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/rozupes-pag/ilhhk.html", "10 €", "Preiļi un raj.", "3 ha.", "76660010146", 46072.54166666667)
)

$row = 2
foreach ($listing in $newListings) {
    $ws1.Range("A" + $lastRow1 + ":F" + $lastRow1).Copy() | Out-Null
    $ws2.Range("A" + $row + ":F" + $row).PasteSpecial(-4122) | Out-Null

    for ($c = 1; $c -le 5; $c++) {
        Set-TextValue $ws2.Cells.Item($row, $c) ([string]$listing[$c - 1])
    }
    $ws2.Cells.Item($row, 6).Value = $listing[5]

    $ws2.Hyperlinks.Add($ws2.Cells.Item($row, 1), [string]$listing[0]) | Out-Null

    $row++
}
